# The commit inserts a new weekly price record as row 12 in the daily
# logic subset, pushing the former rows 12-100 down to 13-101.
#
# Equivalent to Excel's "Insert Copied/Blank Cells" -> Entire Row at row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 12..100 down to 13..101, leaving row 12 blank
# (formatting of the row above/below is carried along automatically).
$ws.Rows(12).Insert()

# Populate the newly inserted row 12 with the new record's data.
$ws.Range("A12").Value = 6
$ws.Range("B12").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C12").Value = "Metropolitana"
$ws.Range("D12").Value = 44473
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = 100112001
$ws.Range("G12").Value = "Berenjena"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 850
$ws.Range("K12").Value = 6000
$ws.Range("L12").Value = 7000
$ws.Range("M12").Value = 6588
$ws.Range("N12").Value = "$/caja 50 unidades"
$ws.Range("O12").Value = "Región de Arica y Parinacota"
$ws.Range("P12").Value = 132
$ws.Range("Q12").Value = 50
$ws.Range("R12").Value = "Hortaliza"
